$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing data rows (19, 20, 21) from the bottom up
$ws.Rows.Item(21).EntireRow.Delete()
$ws.Rows.Item(20).EntireRow.Delete()
$ws.Rows.Item(19).EntireRow.Delete()

$ws.Cells.Item(2,1).Value = 2015
$ws.Cells.Item(2,2).Value = "MS"
$ws.Cells.Item(2,3).Value = "BAY SAINT LOUIS"
$ws.Cells.Item(2,4).Value = 261975
$ws.Cells.Item(2,5).Value = 90096
$ws.Cells.Item(2,6).Value = 97921
$ws.Cells.Item(2,7).Value = -71.23
$ws.Cells.Item(2,8).Value = 320.82
$ws.Cells.Item(3,1).Value = 2011
$ws.Cells.Item(3,2).Value = "FL"
$ws.Cells.Item(3,3).Value = "BRANDON"
$ws.Cells.Item(3,4).Value = 388
$ws.Cells.Item(3,5).Value = 3104
$ws.Cells.Item(3,6).Value = 3598
$ws.Cells.Item(3,7).Value = -57.52
$ws.Cells.Item(3,8).Value = 142.5
$ws.Cells.Item(4,1).Value = 2015
$ws.Cells.Item(4,2).Value = "LA"
$ws.Cells.Item(4,3).Value = "CHARENTON"
$ws.Cells.Item(4,4).Value = 4565
$ws.Cells.Item(4,5).Value = 521
$ws.Cells.Item(4,6).Value = 566
$ws.Cells.Item(4,7).Value = -66.09
$ws.Cells.Item(4,8).Value = 196.43
$ws.Cells.Item(5,1).Value = 1986
$ws.Cells.Item(5,2).Value = "FL"
$ws.Cells.Item(5,3).Value = "DESTIN"
$ws.Cells.Item(5,4).Value = 10707
$ws.Cells.Item(5,5).Value = 7552
$ws.Cells.Item(5,6).Value = 15445
$ws.Cells.Item(5,7).Value = -51.22
$ws.Cells.Item(5,8).Value = 140.34
$ws.Cells.Item(6,1).Value = 1993
$ws.Cells.Item(6,2).Value = "FL"
$ws.Cells.Item(6,3).Value = "GULF BREEZE"
$ws.Cells.Item(6,4).Value = 152279
$ws.Cells.Item(6,5).Value = 192525
$ws.Cells.Item(6,6).Value = 318049
$ws.Cells.Item(6,7).Value = -75.12
$ws.Cells.Item(6,8).Value = 446.8
$ws.Cells.Item(7,1).Value = 1979
$ws.Cells.Item(7,2).Value = "FL"
$ws.Cells.Item(7,3).Value = "ISLAMORADA"
$ws.Cells.Item(7,4).Value = 305953
$ws.Cells.Item(7,5).Value = 473696
$ws.Cells.Item(7,6).Value = 1391085
$ws.Cells.Item(7,7).Value = -94.55
$ws.Cells.Item(7,8).Value = 100.27
$ws.Cells.Item(8,1).Value = 2016
$ws.Cells.Item(8,2).Value = "LA"
$ws.Cells.Item(8,3).Value = "JACKSON"
$ws.Cells.Item(8,4).Value = 1800
$ws.Cells.Item(8,5).Value = 1600
$ws.Cells.Item(8,6).Value = 1722
$ws.Cells.Item(8,7).Value = -75.44
$ws.Cells.Item(8,8).Value = 74.76
$ws.Cells.Item(9,1).Value = 1993
$ws.Cells.Item(9,2).Value = "FL"
$ws.Cells.Item(9,3).Value = "MILTON"
$ws.Cells.Item(9,4).Value = 182920
$ws.Cells.Item(9,5).Value = 833251
$ws.Cells.Item(9,6).Value = 1376521
$ws.Cells.Item(9,7).Value = -64.89
$ws.Cells.Item(9,8).Value = 52.01
$ws.Cells.Item(10,1).Value = 2007
$ws.Cells.Item(10,2).Value = "LA"
$ws.Cells.Item(10,3).Value = "MYRTLE GROVE"
$ws.Cells.Item(10,4).Value = 151
$ws.Cells.Item(10,5).Value = 152
$ws.Cells.Item(10,6).Value = 187
$ws.Cells.Item(10,7).Value = -55.58
$ws.Cells.Item(10,8).Value = 308.11
$ws.Cells.Item(11,1).Value = 2020
$ws.Cells.Item(11,2).Value = "FL"
$ws.Cells.Item(11,3).Value = "PLACIDA"
$ws.Cells.Item(11,4).Value = 1038
$ws.Cells.Item(11,5).Value = 3440
$ws.Cells.Item(11,6).Value = 3440
$ws.Cells.Item(11,7).Value = -50.55
$ws.Cells.Item(11,8).Value = 60.19
$ws.Cells.Item(12,1).Value = 2010
$ws.Cells.Item(12,2).Value = "FL"
$ws.Cells.Item(12,3).Value = "POMONA PARK"
$ws.Cells.Item(12,4).Value = 7681
$ws.Cells.Item(12,5).Value = 10369
$ws.Cells.Item(12,6).Value = 12269
$ws.Cells.Item(12,7).Value = -64.51
$ws.Cells.Item(12,8).Value = 232.08
$ws.Cells.Item(13,1).Value = 2008
$ws.Cells.Item(13,2).Value = "TX"
$ws.Cells.Item(13,3).Value = "PORT MANSFIELD"
$ws.Cells.Item(13,4).Value = 2810
$ws.Cells.Item(13,5).Value = 2248
$ws.Cells.Item(13,6).Value = 2709
$ws.Cells.Item(13,7).Value = -73.76
$ws.Cells.Item(13,8).Value = 85.97
$ws.Cells.Item(14,1).Value = 2010
$ws.Cells.Item(14,2).Value = "AL"
$ws.Cells.Item(14,3).Value = "SEMINOLE"
$ws.Cells.Item(14,4).Value = 4807
$ws.Cells.Item(14,5).Value = 4351
$ws.Cells.Item(14,6).Value = 5148
$ws.Cells.Item(14,7).Value = -78.84
$ws.Cells.Item(14,8).Value = 141.44
$ws.Cells.Item(15,1).Value = 2016
$ws.Cells.Item(15,2).Value = "FL"
$ws.Cells.Item(15,3).Value = "SPRING HILL"
$ws.Cells.Item(15,4).Value = 1073.5
$ws.Cells.Item(15,5).Value = 2161.03
$ws.Cells.Item(15,6).Value = 2325
$ws.Cells.Item(15,7).Value = -51.72
$ws.Cells.Item(15,8).Value = 112.15
$ws.Cells.Item(16,1).Value = 2015
$ws.Cells.Item(16,2).Value = "LA"
$ws.Cells.Item(16,3).Value = "SUGARTOWN"
$ws.Cells.Item(16,4).Value = 2661
$ws.Cells.Item(16,5).Value = 705
$ws.Cells.Item(16,6).Value = 766
$ws.Cells.Item(16,7).Value = -72.71
$ws.Cells.Item(16,8).Value = 151.99
$ws.Cells.Item(17,1).Value = 1987
$ws.Cells.Item(17,2).Value = "FL"
$ws.Cells.Item(17,3).Value = "TAVERNIER"
$ws.Cells.Item(17,4).Value = 80053
$ws.Cells.Item(17,5).Value = 123
$ws.Cells.Item(17,6).Value = 246
$ws.Cells.Item(17,7).Value = -88.9
$ws.Cells.Item(17,8).Value = 637.41
$ws.Cells.Item(18,1).Value = 2000
$ws.Cells.Item(18,2).Value = "LA"
$ws.Cells.Item(18,3).Value = "YOUNGSVILLE"
$ws.Cells.Item(18,4).Value = 1096
$ws.Cells.Item(18,5).Value = 113
$ws.Cells.Item(18,6).Value = 165
$ws.Cells.Item(18,7).Value = -95.8
$ws.Cells.Item(18,8).Value = 149.66
